$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.786.82'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.494.64'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.33'
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.99'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.565'
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.517.74'
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("E10").Value = '  +2.76%  '
$ws.Range("E11").Value = '  -1.82%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.347'
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.944.43'
$ws.Range("E14").Value = '  +1.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.96'
$ws.Range("E15").Value = '  +2.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.677.14'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.521.12'
$ws.Range("E18").Value = '  +2.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.03'
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("E21").Value = '  +1.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.95'
$ws.Range("E23").Value = '  +5.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.93'
$ws.Range("E24").Value = '  +4.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.420'
$ws.Range("E25").Value = '  +4.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.50'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.63'
$ws.Range("E29").Value = '  +2.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0765'
$ws.Range("E30").Value = '  +2.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.75'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '170.68'
$ws.Range("E32").Value = '  +4.66%  '
$ws.Range("E33").Value = '  +10.28%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.37'
$ws.Range("E35").Value = '  +2.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.34'
$ws.Range("E36").Value = '  +1.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.05'
$ws.Range("E37").Value = '  +1.61%  '
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.86'
$ws.Range("E39").Value = '  +1.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.58'
$ws.Range("E41").Value = '  +1.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '283.42'
$ws.Range("E42").Value = '  +4.73%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.20'
$ws.Range("E43").Value = '  +5.43%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.995'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("E45").Value = '  +4.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '130.15'
$ws.Range("E46").Value = '  +8.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.88'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0922'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0502'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0218'
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.33'
$ws.Range("E51").Value = '  +2.61%  '
